$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 29447.5
$ws.Range("J87").Value = 29447.5
$ws.Range("L87").Value = 29447.5
$ws.Range("N87").Value = -31943.5
$ws.Range("H90").Value = 29447.5
$ws.Range("J90").Value = 29447.5
$ws.Range("L90").Value = 88342.5
$ws.Range("N90").Value = -100822.5
$ws.Range("H94").Value = 2300
$ws.Range("I94").Value = 2300
$ws.Range("K94").Value = 2300
$ws.Range("M94").Value = -1849
$ws.Range("H98").Value = 1766.3462
$ws.Range("I98").Value = 1901.1111
$ws.Range("J98").Value = 1463.125
$ws.Range("K98").Value = 1901.1111
$ws.Range("L98").Value = 1463.125
$ws.Range("M98").Value = -403.1111000000001
$ws.Range("N98").Value = -4459.125
$ws.Range("H122").Value = 1766.3462
$ws.Range("I122").Value = 1901.1111
$ws.Range("J122").Value = 1463.125
$ws.Range("K122").Value = 5703.3333
$ws.Range("L122").Value = 4389.375
$ws.Range("M122").Value = -3253.3333
$ws.Range("N122").Value = -9289.375
$ws.Range("H132").Value = 1363725.2
$ws.Range("I132").Value = 2830.4
$ws.Range("J132").Value = 8168199.5
$ws.Range("K132").Value = 8491.200000000001
$ws.Range("L132").Value = 24504598.5
$ws.Range("M132").Value = -5961.200000000001
$ws.Range("N132").Value = -24509658.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2996.9285
$ws.Range("I45").Value = 3036.1943
$ws.Range("J45").Value = 2761.3333
$ws.Range("K45").Value = 3036.1943
$ws.Range("L45").Value = 2761.3333
$ws.Range("M45").Value = -2659.1943
$ws.Range("N45").Value = -3515.3333
$ws.Range("H132").Value = 40159.67
$ws.Range("I132").Value = 27116.334
$ws.Range("J132").Value = 79289.69500000001
$ws.Range("K132").Value = 81349.00199999999
$ws.Range("L132").Value = 237869.085
$ws.Range("M132").Value = -78819.00199999999
$ws.Range("N132").Value = -242929.085

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 29590066
$ws.Range("I105").Value = 33535052
$ws.Range("J105").Value = 2666.5
$ws.Range("K105").Value = 33535052
$ws.Range("L105").Value = 2666.5
$ws.Range("M105").Value = -33533305
$ws.Range("N105").Value = -6160.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 773.1799999999999
$ws.Range("I31").Value = 918.2941
$ws.Range("J31").Value = 743.4578
$ws.Range("K31").Value = 918.2941
$ws.Range("L31").Value = 743.4578
$ws.Range("M31").Value = -623.2941
$ws.Range("N31").Value = -1333.4578
$ws.Range("H34").Value = 773.1799999999999
$ws.Range("I34").Value = 918.2941
$ws.Range("J34").Value = 743.4578
$ws.Range("K34").Value = 918.2941
$ws.Range("L34").Value = 743.4578
$ws.Range("M34").Value = -716.2941
$ws.Range("N34").Value = -1147.4578
$ws.Range("H64").Value = 29692.75
$ws.Range("J64").Value = 29692.75
$ws.Range("L64").Value = 29692.75
$ws.Range("N64").Value = -30188.75
$ws.Range("H67").Value = 29692.75
$ws.Range("J67").Value = 29692.75
$ws.Range("L67").Value = 29692.75
$ws.Range("N67").Value = -31408.75
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 42403.4
$ws.Range("I132").Value = 2345.7222
$ws.Range("J132").Value = 145408.86
$ws.Range("K132").Value = 7037.1666
$ws.Range("L132").Value = 436226.58
$ws.Range("M132").Value = -4507.1666
$ws.Range("N132").Value = -441286.58
$ws.Range("H134").Value = 38067.668
$ws.Range("I134").Value = 1976.5
$ws.Range("J134").Value = 110250
$ws.Range("K134").Value = 5929.5
$ws.Range("L134").Value = 330750
$ws.Range("M134").Value = -3394.5
$ws.Range("N134").Value = -335820

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1076.7462
$ws.Range("I68").Value = 522.44446
$ws.Range("J68").Value = 1720.4517
$ws.Range("K68").Value = 1567.33338
$ws.Range("L68").Value = 5161.355100000001
$ws.Range("M68").Value = -756.33338
$ws.Range("N68").Value = -6783.355100000001
$ws.Range("H71").Value = 1076.7462
$ws.Range("I71").Value = 522.44446
$ws.Range("J71").Value = 1720.4517
$ws.Range("K71").Value = 4702.00014
$ws.Range("L71").Value = 15484.0653
$ws.Range("M71").Value = -646.0001400000001
$ws.Range("N71").Value = -23596.0653
$ws.Range("H107").Value = 956.97144
$ws.Range("I107").Value = 438.5
$ws.Range("J107").Value = 2088.182
$ws.Range("K107").Value = 1315.5
$ws.Range("L107").Value = 6264.545999999999
$ws.Range("M107").Value = 604.5
$ws.Range("N107").Value = -10104.546
$ws.Range("H113").Value = 550.04083
$ws.Range("I113").Value = 499.88235
$ws.Range("J113").Value = 576.6875
$ws.Range("K113").Value = 1499.64705
$ws.Range("L113").Value = 1730.0625
$ws.Range("M113").Value = 670.35295
$ws.Range("N113").Value = -6070.0625
$ws.Range("H125").Value = 3750
$ws.Range("I125").Value = 1333.3334
$ws.Range("J125").Value = 4555.5557
$ws.Range("K125").Value = 4000.0002
$ws.Range("L125").Value = 13666.6671
$ws.Range("M125").Value = 919.9998000000001
$ws.Range("N125").Value = -23506.6671
$ws.Range("H136").Value = 2343.3333
$ws.Range("I136").Value = 2343.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7029.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1929.999899999999
$ws.Range("N136").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 17441.5
$ws.Range("J95").Value = 17441.5
$ws.Range("L95").Value = 17441.5
$ws.Range("N95").Value = -22933.5
$ws.Range("H102").Value = 1776.3572
$ws.Range("I102").Value = 1731.8948
$ws.Range("J102").Value = 1870.2222
$ws.Range("K102").Value = 1731.8948
$ws.Range("L102").Value = 1870.2222
$ws.Range("M102").Value = -109.8948
$ws.Range("N102").Value = -5114.2222
$ws.Range("H122").Value = 2500.3572
$ws.Range("I122").Value = 2188.5881
$ws.Range("K122").Value = 6565.7643
$ws.Range("M122").Value = -4115.7643
$ws.Range("H132").Value = 61536.793
$ws.Range("I132").Value = 45767.39
$ws.Range("J132").Value = 94509.17999999999
$ws.Range("K132").Value = 137302.17
$ws.Range("L132").Value = 283527.54
$ws.Range("M132").Value = -134772.17
$ws.Range("N132").Value = -288587.54
$ws.Range("H134").Value = 21950.285
$ws.Range("I134").Value = 29000
$ws.Range("J134").Value = 20775.334
$ws.Range("K134").Value = 87000
$ws.Range("L134").Value = 62326.00199999999
$ws.Range("M134").Value = -84465
$ws.Range("N134").Value = -67396.00199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2585.5715
$ws.Range("I7").Value = 2659.8
$ws.Range("K7").Value = 2659.8
$ws.Range("M7").Value = -2547.8
$ws.Range("H40").Value = 3883.5625
$ws.Range("I40").Value = 4777.9
$ws.Range("J40").Value = 2393
$ws.Range("K40").Value = 4777.9
$ws.Range("L40").Value = 2393
$ws.Range("M40").Value = -4641.9
$ws.Range("N40").Value = -2665
$ws.Range("H126").Value = 2585.5715
$ws.Range("I126").Value = 2659.8
$ws.Range("K126").Value = 7979.400000000001
$ws.Range("M126").Value = -5509.400000000001
$ws.Range("H132").Value = 41376
$ws.Range("I132").Value = 1838.9231
$ws.Range("J132").Value = 80913.08
$ws.Range("K132").Value = 5516.7693
$ws.Range("L132").Value = 242739.24
$ws.Range("M132").Value = -2986.7693
$ws.Range("N132").Value = -247799.24
$ws.Range("H136").Value = 87999.44
$ws.Range("I136").Value = 42391.24
$ws.Range("J136").Value = 214688.89
$ws.Range("K136").Value = 127173.72
$ws.Range("L136").Value = 644066.67
$ws.Range("M136").Value = -124623.72
$ws.Range("N136").Value = -649166.67

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 39000
$ws.Range("J103").Value = 39000
$ws.Range("L103").Value = 39000
$ws.Range("N103").Value = -41344
$ws.Range("H132").Value = 59267.887
$ws.Range("I132").Value = 43719.332
$ws.Range("J132").Value = 93192
$ws.Range("K132").Value = 131157.996
$ws.Range("L132").Value = 279576
$ws.Range("M132").Value = -128627.996
$ws.Range("N132").Value = -284636
